# "Add files via upload" - refresh the Station Standard data:
#  - Rows 20/21 (REFILL STATION B19/B20): Config/Type change from
#    "Atlas Box & Bond Bags" to "Skyshop" / "Duty Free" respectively, and the
#    shared "Drawer Avg" formula is overwritten with its literal result.
#  - Rows 34/35/46/47 (REFILL STATION C33/C34/D45/D46): Config/Type change
#    from "Service Cart"/"Service Cart Single" to "Atlas Box & Bond Bags",
#    and the counts are reset to 0.
#  - Row 56 (REFILL STATION D55): counts reset to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Station")

# --- Row 20: REFILL STATION B19 -> Skyshop -----------------------------
$ws.Range("B20").Value = "Skyshop"
$ws.Range("C20").Value = "Skyshop"
$ws.Range("E20").Value = 14

# --- Row 21: REFILL STATION B20 -> Duty Free ---------------------------
$ws.Range("B21").Value = "Duty Free"
$ws.Range("C21").Value = "Duty Free"
$ws.Range("E21").Value = 5

# --- Row 34: REFILL STATION C33 -> Atlas Box & Bond Bags, zeroed -------
$ws.Range("B34").Value = "Atlas Box & Bond Bags"
$ws.Range("C34").Value = "Atlas Box & Bond Bags"
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0

# --- Row 35: REFILL STATION C34 -> Atlas Box & Bond Bags, zeroed -------
$ws.Range("B35").Value = "Atlas Box & Bond Bags"
$ws.Range("C35").Value = "Atlas Box & Bond Bags"
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0

# --- Row 46: REFILL STATION D45 -> Atlas Box & Bond Bags, zeroed -------
$ws.Range("B46").Value = "Atlas Box & Bond Bags"
$ws.Range("C46").Value = "Atlas Box & Bond Bags"
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0

# --- Row 47: REFILL STATION D46 -> Atlas Box & Bond Bags, zeroed -------
$ws.Range("B47").Value = "Atlas Box & Bond Bags"
$ws.Range("C47").Value = "Atlas Box & Bond Bags"
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0

# --- Row 56: REFILL STATION D55, counts zeroed (Config/Type unchanged) -
$ws.Range("E56").Value = 0
$ws.Range("F56").Value = 0

# --- Refresh the view: scroll so row 7 is at the top, select E1 --------
$ws.Activate()
$excel.Goto($ws.Range("A7"), $true) | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E1").Select() | Out-Null

Write-Host "Station Standard refreshed"
